$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C33").Value = "1 - 4 Only"
$ws.Range("C34").Value = "1 - 4 Only"
$ws.Range("C35").Value = "Holiday Saturday"
$ws.Range("C36").Value = "Overriding and Final Keyword with Variable, mathod and class"
